# Weekly update: two new daily price records were reported for the
# "Vega Modelo de Temuco" / Kiwi sheet. They land between the existing
# row 229 (Hayward / Primera, 2021-08-06) and the old row 230
# (Hayward / Especial, 2021-08-05), so two new rows are inserted at
# position 230, pushing all the following rows down by two and
# extending the used range from A1:T311 to A1:T313.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right before the current row 230. Excel shifts
# rows 230:311 down to 232:313 and carries the date-column (D) number
# format down with them.
$ws.Rows("230:231").Insert()

# ---- New row 230 ----------------------------------------------------
$ws.Range("A230").Value = 10
$ws.Range("B230").Value = "Vega Modelo de Temuco"
$ws.Range("C230").Value = "La Araucanía"
$ws.Range("D230").Value = 44460
$ws.Range("E230").Value = 9
$ws.Range("F230").Value = "Fruta"
$ws.Range("G230").Value = 100101
$ws.Range("H230").Value = "Berries"
$ws.Range("I230").Value = 100101007
$ws.Range("J230").Value = "Kiwi"
$ws.Range("K230").Value = "Hayward"
$ws.Range("L230").Value = "Extra"
$ws.Range("M230").Value = 90
$ws.Range("N230").Value = 13000
$ws.Range("O230").Value = 14000
$ws.Range("P230").Value = 13444
$ws.Range("Q230").Value = "`$/bandeja 10 kilos"
$ws.Range("R230").Value = "Región de O'Higgins"
$ws.Range("S230").Value = 1344
$ws.Range("T230").Value = 10

# ---- New row 231 ----------------------------------------------------
$ws.Range("A231").Value = 10
$ws.Range("B231").Value = "Vega Modelo de Temuco"
$ws.Range("C231").Value = "La Araucanía"
$ws.Range("D231").Value = 44460
$ws.Range("E231").Value = 9
$ws.Range("F231").Value = "Fruta"
$ws.Range("G231").Value = 100101
$ws.Range("H231").Value = "Berries"
$ws.Range("I231").Value = 100101007
$ws.Range("J231").Value = "Kiwi"
$ws.Range("K231").Value = "Hayward"
$ws.Range("L231").Value = "Extra"
$ws.Range("M231").Value = 50
$ws.Range("N231").Value = 21000
$ws.Range("O231").Value = 22000
$ws.Range("P231").Value = 21600
$ws.Range("Q231").Value = "`$/bandeja 18 kilos"
$ws.Range("R231").Value = "Región de O'Higgins"
$ws.Range("S231").Value = 1200
$ws.Range("T231").Value = 18
